$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. The last paragraph ("今天天气不错,心情也不错") currently carries a
#    paragraph-mark rPr (<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>).
#    Re-insert its own content without that pPr so the mark formatting is dropped.
$lastPar = $d.Paragraphs.Last
$lastRange = $lastPar.Range
$lastXml = '<w:p ' + $wns + '>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今天天气不错</w:t></w:r>' `
  + '<w:r w:rsidR="008D6A6A"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>,心情也不错</w:t></w:r>' `
  + '</w:p>'
[void]$lastRange.InsertXML($lastXml)

# 2. Append two brand-new paragraphs at the very end of the document:
#      "星期五"
#      "晴天，今天学习了分支管理，创建了一个dev分支 使用Git 创建分支简单又快捷"
$endRange = $d.Content
$endRange.Collapse(0)
$newXml = '<w:p ' + $wns + '>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>星期五</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p ' + $wns + '>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>晴天，今天学习了分支管理，创建了一个dev分支</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>使用Git</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>创建分支简单又快捷</w:t></w:r>' `
  + '</w:p>'
[void]$endRange.InsertXML($newXml)
